# Trade #21 (global trade #51, zero-indexed #51 -> "Trade #51" row) closed
# at 2026-02-18 00:12:41 - unknown UNKNOWN +0.000%
#
# This script applies the following to the live trading results workbook:
#   1. Updates the rolled-up Summary metrics.
#   2. Updates the MarketMaking row in the Strategy Status sheet.
#   3. Closes the open MarketMaking trade (Trade #51) on both the
#      "All Trades" sheet and the per-strategy "MarketMaking" sheet.
#   4. Appends two freshly logged OPEN trades (#80 momentum, #81
#      HighProbConvergence) to "All Trades" and their respective
#      per-strategy sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value to a cell while forcing a Text number
# format first, so Excel does not auto-coerce date-looking strings
# (e.g. "2026-02-18") into date serials.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# =======================================================================
# 1. Summary sheet
# =======================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.87
$summary.Range("B4").Value = 0.97
$summary.Range("B5").Value = 0.4
$summary.Range("B6").Value = 49
$summary.Range("B7").Value = 28
$summary.Range("B9").Value = 57.14

# =======================================================================
# 2. Strategy Status sheet - MarketMaking row (row 6)
# =======================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.87
$status.Range("D6").Value = 20
$status.Range("E6").Value = 0.06
$status.Range("F6").Value = -0.13
$status.Range("G6").Value = 60

# =======================================================================
# 3. All Trades sheet - close Trade #51 (row 52) + append new trades
# =======================================================================
$allTrades = $wb.Worksheets.Item("All Trades")

# --- Close existing MarketMaking trade (row 52, Trade # = 51) ---
$allTrades.Range("G52").Value = 0.82
$allTrades.Range("H52").Value = "CLOSED"
$allTrades.Range("I52").Value = 20.5882
$allTrades.Range("J52").Value = 0.14
$allTrades.Range("K52").Value = 99.87
$allTrades.Range("L52").Value = "early_exit"
$allTrades.Range("M52").Value = 0.16

# --- Append Trade #80 (momentum, row 81) ---
$allTrades.Range("A81").Value = 80
Set-TextCell $allTrades 81 2 "2026-02-18"
Set-TextCell $allTrades 81 3 "00:12:35"
$allTrades.Range("D81").Value = "momentum"
$allTrades.Range("E81").Value = "DOWN"
$allTrades.Range("F81").Value = 0.68
$allTrades.Range("H81").Value = "OPEN"
$allTrades.Range("I81").Value = 0
$allTrades.Range("J81").Value = 0
$allTrades.Range("K81").Value = 100
$allTrades.Range("M81").Value = 0
$allTrades.Range("N81").Value = 0
$allTrades.Range("O81").Value = 0
$allTrades.Range("P81").Value = 0.9
$allTrades.Range("Q81").Value = "Downward momentum: -1.980% over 10 samples"

# --- Append Trade #81 (HighProbConvergence, row 82) ---
$allTrades.Range("A82").Value = 81
Set-TextCell $allTrades 82 2 "2026-02-18"
Set-TextCell $allTrades 82 3 "00:12:36"
$allTrades.Range("D82").Value = "HighProbConvergence"
$allTrades.Range("E82").Value = "UP"
$allTrades.Range("F82").Value = 0.31
$allTrades.Range("H82").Value = "OPEN"
$allTrades.Range("I82").Value = 0
$allTrades.Range("J82").Value = 0
$allTrades.Range("K82").Value = 100
$allTrades.Range("M82").Value = 0
$allTrades.Range("N82").Value = 0
$allTrades.Range("O82").Value = 0
$allTrades.Range("P82").Value = 0.95
$allTrades.Range("Q82").Value = "Mean reversion UP: price 1.69% below mean (z=-2.38)"

# =======================================================================
# 4. momentum sheet - append Trade #80 (row 13)
# =======================================================================
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A13").Value = 80
Set-TextCell $momentum 13 2 "2026-02-18"
Set-TextCell $momentum 13 3 "00:12:35"
$momentum.Range("D13").Value = "momentum"
$momentum.Range("E13").Value = "DOWN"
$momentum.Range("F13").Value = 0.68
$momentum.Range("H13").Value = "OPEN"
$momentum.Range("I13").Value = 0
$momentum.Range("J13").Value = 0
$momentum.Range("K13").Value = 100
$momentum.Range("L13").Value = 0
$momentum.Range("M13").Value = 0
$momentum.Range("N13").Value = 0.9
$momentum.Range("O13").Value = "Downward momentum: -1.980% over 10 samples"
$momentum.Range("Q13").Value = 0

# =======================================================================
# 5. HighProbConvergence sheet - append Trade #81 (row 7)
# =======================================================================
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("A7").Value = 81
Set-TextCell $hpc 7 2 "2026-02-18"
Set-TextCell $hpc 7 3 "00:12:36"
$hpc.Range("D7").Value = "HighProbConvergence"
$hpc.Range("E7").Value = "UP"
$hpc.Range("F7").Value = 0.31
$hpc.Range("H7").Value = "OPEN"
$hpc.Range("I7").Value = 0
$hpc.Range("J7").Value = 0
$hpc.Range("K7").Value = 100
$hpc.Range("L7").Value = 0
$hpc.Range("M7").Value = 0
$hpc.Range("N7").Value = 0.95
$hpc.Range("O7").Value = "Mean reversion UP: price 1.69% below mean (z=-2.38)"
$hpc.Range("Q7").Value = 0

# =======================================================================
# 6. MarketMaking sheet - close Trade #51 (row 23)
# =======================================================================
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G23").Value = 0.82
$mm.Range("H23").Value = "CLOSED"
$mm.Range("I23").Value = 20.5882
$mm.Range("J23").Value = 0.14
$mm.Range("K23").Value = 99.87
$mm.Range("P23").Value = "early_exit"
$mm.Range("Q23").Value = 0.16
